# Auto-generated edit script applying the Spriggan_Profits price-refresh diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds leve-profit data;
# a scheduled runner refreshed the market-price columns (H-N) for a handful of rows.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 41
$ws.Range("H41").Value = 1045.3077
$ws.Range("I41").Value = 1381.3334
$ws.Range("K41").Value = 1381.3334
$ws.Range("M41").Value = -941.3334
# row 70
$ws.Range("H70").Value = 7243.75
$ws.Range("I70").Value = 1700
$ws.Range("J70").Value = 8035.7144
$ws.Range("K70").Value = 5100
$ws.Range("L70").Value = 24107.1432
$ws.Range("M70").Value = -4830
$ws.Range("N70").Value = -24647.1432
# row 73
$ws.Range("H73").Value = 7243.75
$ws.Range("I73").Value = 1700
$ws.Range("J73").Value = 8035.7144
$ws.Range("K73").Value = 5100
$ws.Range("L73").Value = 24107.1432
$ws.Range("M73").Value = -4164
$ws.Range("N73").Value = -25979.1432
# row 76
$ws.Range("H76").Value = 3994.5
$ws.Range("I76").Value = 3994.5
$ws.Range("K76").Value = 3994.5
$ws.Range("M76").Value = -3679.5
# row 79
$ws.Range("H79").Value = 3994.5
$ws.Range("I79").Value = 3994.5
$ws.Range("K79").Value = 3994.5
$ws.Range("M79").Value = -2902.5
# row 96
$ws.Range("H96").Value = 471.8
$ws.Range("I96").Value = 582.3333
$ws.Range("J96").Value = 306
$ws.Range("K96").Value = 1746.9999
$ws.Range("L96").Value = 918
$ws.Range("M96").Value = -373.9999
$ws.Range("N96").Value = -3664
# row 107
$ws.Range("H107").Value = 1517.4445
$ws.Range("I107").Value = 1058.25
$ws.Range("K107").Value = 1058.25
$ws.Range("M107").Value = 861.75
# row 129
$ws.Range("H129").Value = 2536.3
$ws.Range("I129").Value = 2688.5
$ws.Range("J129").Value = 2308
$ws.Range("K129").Value = 8065.5
$ws.Range("L129").Value = 6924
$ws.Range("M129").Value = -3065.5
$ws.Range("N129").Value = -16924
# row 132
$ws.Range("H132").Value = 2130.9768
$ws.Range("I132").Value = 2169.8096
$ws.Range("K132").Value = 6509.4288
$ws.Range("M132").Value = -3979.4288
# row 138
$ws.Range("H138").Value = 3798.7334
$ws.Range("I138").Value = 2490.875
$ws.Range("J138").Value = 4274.3184
$ws.Range("K138").Value = 7472.625
$ws.Range("L138").Value = 12822.9552
$ws.Range("M138").Value = -2332.625
$ws.Range("N138").Value = -23102.9552

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 43
$ws.Range("H43").Value = 61598
$ws.Range("J43").Value = 88996.5
$ws.Range("L43").Value = 88996.5
$ws.Range("N43").Value = -89622.5
# row 61
$ws.Range("H61").Value = 76928120
$ws.Range("J61").Value = 9197.6
$ws.Range("L61").Value = 9197.6
$ws.Range("N61").Value = -9621.6
# row 74
$ws.Range("H74").Value = 25001710
$ws.Range("I74").Value = 32259668
$ws.Range("K74").Value = 32259668
$ws.Range("M74").Value = -32258794
# row 77
$ws.Range("H77").Value = 25001710
$ws.Range("I77").Value = 32259668
$ws.Range("K77").Value = 161298340
$ws.Range("M77").Value = -161293972
# row 97
$ws.Range("H97").Value = 447.77777
$ws.Range("I97").Value = 447.77777
$ws.Range("K97").Value = 447.77777
$ws.Range("M97").Value = 48.22223000000002
# row 110
$ws.Range("H110").Value = 50990.1
$ws.Range("I110").Value = 56570.723
$ws.Range("K110").Value = 56570.723
$ws.Range("M110").Value = -54525.723
# row 122
$ws.Range("H122").Value = 5308.625
$ws.Range("I122").Value = 2432.1667
$ws.Range("J122").Value = 13938
$ws.Range("K122").Value = 7296.500100000001
$ws.Range("L122").Value = 41814
$ws.Range("M122").Value = -4846.500100000001
$ws.Range("N122").Value = -46714
# row 132
$ws.Range("H132").Value = 2633734.8
$ws.Range("J132").Value = 2512.5715
$ws.Range("L132").Value = 7537.7145
$ws.Range("N132").Value = -12597.7145
# row 136
$ws.Range("H136").Value = 76928120
$ws.Range("J136").Value = 9197.6
$ws.Range("L136").Value = 27592.8
$ws.Range("N136").Value = -32692.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 3751.3333
$ws.Range("I86").Value = 3683.2727
$ws.Range("K86").Value = 3683.2727
$ws.Range("M86").Value = -2560.2727
# row 89
$ws.Range("H89").Value = 3751.3333
$ws.Range("I89").Value = 3683.2727
$ws.Range("K89").Value = 18416.3635
$ws.Range("M89").Value = -12800.3635
# row 103
$ws.Range("H103").Value = 45108.332
$ws.Range("J103").Value = 45108.332
$ws.Range("L103").Value = 45108.332
$ws.Range("N103").Value = -47452.332
# row 105
$ws.Range("H105").Value = 3763.625
$ws.Range("I105").Value = 3552.25
$ws.Range("K105").Value = 3552.25
$ws.Range("M105").Value = -1805.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 7921.231
$ws.Range("I7").Value = 10257.6
$ws.Range("K7").Value = 10257.6
$ws.Range("M7").Value = -10144.6
# row 16
$ws.Range("H16").Value = 18365.5
$ws.Range("I16").Value = 25749
$ws.Range("K16").Value = 25749
$ws.Range("M16").Value = -25462
# row 22
$ws.Range("H22").Value = 494.25
$ws.Range("I22").Value = 494.25
$ws.Range("K22").Value = 494.25
$ws.Range("M22").Value = -144.25
# row 31
$ws.Range("H31").Value = 6072.7856
$ws.Range("I31").Value = 8803.777
$ws.Range("J31").Value = 1157
$ws.Range("K31").Value = 8803.777
$ws.Range("L31").Value = 1157
$ws.Range("M31").Value = -8508.777
$ws.Range("N31").Value = -1747
# row 34
$ws.Range("H34").Value = 6072.7856
$ws.Range("I34").Value = 8803.777
$ws.Range("J34").Value = 1157
$ws.Range("K34").Value = 8803.777
$ws.Range("L34").Value = 1157
$ws.Range("M34").Value = -8601.777
$ws.Range("N34").Value = -1561
# row 107
$ws.Range("H107").Value = 39122.652
$ws.Range("I107").Value = 258.0625
$ws.Range("K107").Value = 258.0625
$ws.Range("M107").Value = 1661.9375
# row 113
$ws.Range("H113").Value = 18365.5
$ws.Range("I113").Value = 25749
$ws.Range("K113").Value = 25749
$ws.Range("M113").Value = -23579

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 133
$ws.Range("H133").Value = 15708.25
$ws.Range("J133").Value = 18277.666
$ws.Range("L133").Value = 54832.99800000001
$ws.Range("N133").Value = -64952.99800000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 19
$ws.Range("H19").Value = 9551.5
$ws.Range("I19").Value = 9702
$ws.Range("K19").Value = 9702
$ws.Range("M19").Value = -9414
# row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# row 80
$ws.Range("H80").Value = 865
$ws.Range("J80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("N80").Value = -2996
# row 83
$ws.Range("H83").Value = 865
$ws.Range("J83").Value = 1000
$ws.Range("L83").Value = 5000
$ws.Range("N83").Value = -14984
# row 107
$ws.Range("H107").Value = 1433.5
$ws.Range("I107").Value = 1141.6111
$ws.Range("J107").Value = 2309.1667
$ws.Range("K107").Value = 1141.6111
$ws.Range("L107").Value = 2309.1667
$ws.Range("M107").Value = 778.3888999999999
$ws.Range("N107").Value = -6149.1667
# row 113
$ws.Range("H113").Value = 92264.27
$ws.Range("I113").Value = 101290.7
$ws.Range("K113").Value = 101290.7
$ws.Range("M113").Value = -99120.7
# row 132
$ws.Range("H132").Value = 9617749
$ws.Range("I132").Value = 9617749
$ws.Range("K132").Value = 28853247
$ws.Range("M132").Value = -28850717

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 4662.5
$ws.Range("I61").Value = 4662.5
$ws.Range("K61").Value = 4662.5
$ws.Range("M61").Value = -4460.5
# row 70
$ws.Range("H70").Value = 44444
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 44444
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 44444
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -44984
# row 73
$ws.Range("H73").Value = 44444
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 44444
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 44444
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -46316
# row 74
$ws.Range("H74").Value = 55554.8
$ws.Range("I74").Value = 61887.5
$ws.Range("J74").Value = 51333
$ws.Range("K74").Value = 61887.5
$ws.Range("L74").Value = 51333
$ws.Range("M74").Value = -60889.5
$ws.Range("N74").Value = -53329
# row 77
$ws.Range("H77").Value = 55554.8
$ws.Range("I77").Value = 61887.5
$ws.Range("J77").Value = 51333
$ws.Range("K77").Value = 185662.5
$ws.Range("L77").Value = 153999
$ws.Range("M77").Value = -180670.5
$ws.Range("N77").Value = -163983
# row 113
$ws.Range("H113").Value = 4662.5
$ws.Range("I113").Value = 4662.5
$ws.Range("K113").Value = 4662.5
$ws.Range("M113").Value = -2492.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 97
$ws.Range("H97").Value = 28786
$ws.Range("J97").Value = 28786
$ws.Range("L97").Value = 28786
$ws.Range("N97").Value = -30768
# row 103
$ws.Range("H103").Value = 40050.75
$ws.Range("J103").Value = 40050.75
$ws.Range("L103").Value = 40050.75
$ws.Range("N103").Value = -42394.75
# row 132
$ws.Range("H132").Value = 71429610
$ws.Range("I132").Value = 83334450
$ws.Range("K132").Value = 250003350
$ws.Range("M132").Value = -250000820
